# Update Ankit Rajpoot's per-innings batting activity figures.
# Row 2 "runs" (C2) and "sixes" (F2) swap with row 3's "runs" (C3) and
# "sixes" (F3) values. The source sheet stores these numeric-looking
# figures as text, so each value is entered with a leading apostrophe to
# keep it text (matching how it was authored originally) rather than
# letting Excel auto-convert it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'2"
$ws.Range("F2").Value = "'0"

$ws.Range("C3").Value = "'7"
$ws.Range("F3").Value = "'1"
